$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("DataFetchFlag")
$ws2 = $wb.Worksheets.Item("DataFetchXL")

# New test rows: IssuedForReview / RequestForInformation / IssuedForApproval
$names = @(
    "FLD_Transmittals_New_IssuedForReview",
    "FLD_Transmittals_New_RequestForInformation",
    "FLD_Transmittals_New_IssuedForApproval"
)

# Path shown in the cell (display text) - no "file:///" prefix, matches the
# existing rows' convention.
$paths = @(
    "\\src\com\proj\suiteTRANSMITTALS\testdata\TransmittalsTestData-IssuedForReview.xlsx",
    "\\src\com\proj\suiteTRANSMITTALS\testdata\TransmittalsTestData-RequestForInformation.xlsx",
    "\\src\com\proj\suiteTRANSMITTALS\testdata\TransmittalsTestData-IssuedForApproval.xlsx"
)

# Full hyperlink target (what the link actually navigates to).
$links = @(
    "file:///\\src\com\proj\suiteTRANSMITTALS\testdata\TransmittalsTestData-IssuedForReview.xlsx",
    "file:///\\src\com\proj\suiteTRANSMITTALS\testdata\TransmittalsTestData-RequestForInformation.xlsx",
    "file:///\\src\com\proj\suiteTRANSMITTALS\testdata\TransmittalsTestData-IssuedForApproval.xlsx"
)

# Write column A (TestCaseName) on both sheets first, then column B/C, so the
# shared-string table fills up in the same order the rows are authored in.
for ($i = 0; $i -lt 3; $i++) {
    $row = 8 + $i
    $ws1.Range("A$row").Value = $names[$i]
    $ws2.Range("A$row").Value = $names[$i]
}

for ($i = 0; $i -lt 3; $i++) {
    $row = 8 + $i

    # Sheet "DataFetchFlag": DataFetchFlag column is always "XL" here.
    $ws1.Range("B$row").Value = "XL"

    # Sheet "DataFetchXL": ExcelDataSheetPath + FirstSheetName columns.
    $ws2.Range("B$row").Value = $paths[$i]
    $ws2.Range("C$row").Value = "Transmittals_New"

    # Wire the hyperlink up, then restore the worksheet's Hyperlink cell
    # style so the new cell matches the look of the existing linked cells.
    $ws2.Hyperlinks.Add($ws2.Range("B$row"), $links[$i])
    $ws2.Range("B$row").Style = "Hyperlink"
}

$excel.CutCopyMode = 0

# Keep the data-validation list in sync with the newly added rows.
$ws1.Range("B2:B10").Validation.Delete()
$ws1.Range("B2:B10").Validation.Add(3, 1, 1, """XL,DB""")
$ws1.Range("B2:B10").Validation.IgnoreBlank = 1
$ws1.Range("B2:B10").Validation.InCellDropdown = 1
$ws1.Range("B2:B10").Validation.ShowInput = 1
$ws1.Range("B2:B10").Validation.ShowError = 1
